# Update the "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages rebuild at 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 158
    $ws.Range("F3").Value = 1732
    $ws.Range("F7").Value = 12033
    $ws.Range("F14").Value = 13495
    $ws.Range("F15").Value = 13523
    $ws.Range("F23").Value = 1921
    $ws.Range("F24").Value = 179
}

# F20 diverges slightly between the two sheets in the regenerated output.
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F20").Value = 972

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F20").Value = 974
